# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns retain their text formatting,
# since the source data are inline strings (not numeric values).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.725.81'
$ws.Cells.Item(2, 5).Value = '  +4.63%  '
$ws.Cells.Item(3, 4).Value = '1.873.86'
$ws.Cells.Item(3, 5).Value = '  +3.64%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).Value = '283.02'
$ws.Cells.Item(5, 5).Value = '  +2.28%  '
$ws.Cells.Item(6, 5).Value = '  -0.17%  '
$ws.Cells.Item(7, 4).Value = '0.5191'
$ws.Cells.Item(7, 5).Value = '  +3.54%  '
$ws.Cells.Item(8, 4).Value = '0.3529'
$ws.Cells.Item(8, 5).Value = '  +0.73%  '
$ws.Cells.Item(9, 4).Value = '45.28'
$ws.Cells.Item(9, 5).Value = '  +3.29%  '
$ws.Cells.Item(10, 4).Value = '0.07100'
$ws.Cells.Item(10, 5).Value = '  +6.63%  '
$ws.Cells.Item(11, 4).Value = '20.24'
$ws.Cells.Item(11, 5).Value = '  +1.41%  '
$ws.Cells.Item(12, 4).Value = '0.8214'
$ws.Cells.Item(12, 5).Value = '  -1.83%  '
$ws.Cells.Item(13, 4).Value = '0.07762'
$ws.Cells.Item(13, 5).Value = '  -0.96%  '
$ws.Cells.Item(14, 4).Value = '1.863.50'
$ws.Cells.Item(14, 5).Value = '  +3.15%  '
$ws.Cells.Item(15, 4).Value = '5.171'
$ws.Cells.Item(15, 5).Value = '  +2.46%  '
$ws.Cells.Item(16, 5).Value = '  +3.03%  '
$ws.Cells.Item(17, 5).Value = '  -0.10%  '
$ws.Cells.Item(18, 5).Value = '  +4.13%  '
$ws.Cells.Item(19, 4).Value = '0.000008163'
$ws.Cells.Item(19, 5).Value = '  +3.56%  '
$ws.Cells.Item(20, 4).Value = '0.9994'
$ws.Cells.Item(20, 5).Value = '  -0.17%  '
$ws.Cells.Item(21, 4).Value = '26.770.70'
$ws.Cells.Item(21, 5).Value = '  +4.50%  '
$ws.Cells.Item(22, 4).Value = '4.792'
$ws.Cells.Item(22, 5).Value = '  +1.64%  '
$ws.Cells.Item(23, 5).Value = '  +2.10%  '
$ws.Cells.Item(24, 4).Value = '6.242'
$ws.Cells.Item(24, 5).Value = '  +3.05%  '
$ws.Cells.Item(25, 4).Value = '2.428'
$ws.Cells.Item(25, 5).Value = '  +15.10%  '
$ws.Cells.Item(26, 4).Value = '145.74'
$ws.Cells.Item(26, 5).Value = '  +3.02%  '
$ws.Cells.Item(27, 5).Value = '  +3.17%  '
$ws.Cells.Item(28, 4).Value = '1.669'
$ws.Cells.Item(28, 5).Value = '  +0.44%  '
$ws.Cells.Item(29, 4).Value = '111.46'
$ws.Cells.Item(29, 5).Value = '  +2.58%  '
$ws.Cells.Item(30, 4).Value = '4.423'
$ws.Cells.Item(31, 4).Value = '4.362'
$ws.Cells.Item(31, 5).Value = '  +3.79%  '
$ws.Cells.Item(32, 4).Value = '0.08857'
$ws.Cells.Item(32, 5).Value = '  +0.36%  '
$ws.Cells.Item(33, 4).Value = '0.04916'
$ws.Cells.Item(33, 5).Value = '  +2.38%  '
$ws.Cells.Item(34, 4).Value = '1.179'
$ws.Cells.Item(34, 5).Value = '  +5.11%  '
$ws.Cells.Item(35, 4).Value = '0.7470'
$ws.Cells.Item(35, 5).Value = '  +1.89%  '
$ws.Cells.Item(36, 4).Value = '3.301'
$ws.Cells.Item(36, 5).Value = '  +8.86%  '
$ws.Cells.Item(37, 4).Value = '2.865'
$ws.Cells.Item(37, 5).Value = '  +0.58%  '
$ws.Cells.Item(38, 4).Value = '2.418'
$ws.Cells.Item(38, 5).Value = '  +4.52%  '
$ws.Cells.Item(39, 4).Value = '0.5325'
$ws.Cells.Item(39, 5).Value = '  +2.43%  '
$ws.Cells.Item(40, 4).Value = '0.01883'
$ws.Cells.Item(40, 5).Value = '  +1.24%  '
$ws.Cells.Item(41, 4).Value = '0.9756'
$ws.Cells.Item(41, 5).Value = '  +1.93%  '
$ws.Cells.Item(42, 4).Value = '116.46'
$ws.Cells.Item(42, 5).Value = '  +3.88%  '
$ws.Cells.Item(43, 4).Value = '6.310'
$ws.Cells.Item(43, 5).Value = '  +2.34%  '
$ws.Cells.Item(44, 5).Value = '  +2.04%  '
$ws.Cells.Item(45, 4).Value = '0.9994'
$ws.Cells.Item(45, 5).Value = '  -0.18%  '
$ws.Cells.Item(46, 4).Value = '0.4622'
$ws.Cells.Item(46, 5).Value = '  +0.78%  '
$ws.Cells.Item(47, 4).Value = '0.1370'
$ws.Cells.Item(47, 5).Value = '  -0.69%  '
$ws.Cells.Item(48, 4).Value = '9.489'
$ws.Cells.Item(48, 5).Value = '  +3.08%  '
$ws.Cells.Item(49, 4).Value = '36.71'
$ws.Cells.Item(49, 5).Value = '  +2.98%  '
$ws.Cells.Item(50, 5).Value = '  +1.99%  '
$ws.Cells.Item(51, 4).Value = '0.05933'
$ws.Cells.Item(51, 5).Value = '  +1.66%  '
